$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns I and J
# Copy the formatting from the existing header cell (H1) so the new
# header cells share the same bold/centered/bordered style.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-26 for columns I (I0) and J (IF)
$data = @{
    2  = @(7, 8)
    3  = @(6, 7)
    4  = @(6, 6)
    5  = @(6, 6)
    6  = @(7, 7)
    7  = @(4, 4)
    8  = @(8, 9)
    9  = @(9, 9)
    10 = @(2, 3)
    11 = @(8, 8)
    12 = @(7, 7)
    13 = @(7, 7)
    14 = @(7, 8)
    15 = @(8, 8)
    16 = @(7, 7)
    17 = @(8, 9)
    18 = @(8, 8)
    19 = @(6, 6)
    20 = @(8, 8)
    21 = @(6, 6)
    22 = @(9, 9)
    23 = @(7, 7)
    24 = @(8, 9)
    25 = @(6, 6)
    26 = @(5, 5)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
